# Weekly fruit/vegetable data update.
# Two new daily records were inserted ahead of the existing history for this
# market/product series, pushing every subsequent row down by two positions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 217-218; everything below (old rows 217-313)
# shifts down to 219-315, carrying its formatting (incl. the date style on
# column D) along with it.
$ws.Rows("217:218").Insert()

# New row 217
$ws.Range("A217").Value = 9
$ws.Range("B217").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C217").Value = 'Metropolitana'
$ws.Range("D217").Value = 44845
$ws.Range("E217").Value = 13
$ws.Range("F217").Value = 100112001
$ws.Range("G217").Value = 'Berenjena'
$ws.Range("H217").Value = 'Sin especificar'
$ws.Range("I217").Value = 'Primera'
$ws.Range("J217").Value = 180
$ws.Range("K217").Value = 12000
$ws.Range("L217").Value = 14000
$ws.Range("M217").Value = 13111
$ws.Range("N217").Value = '$/caja 50 unidades'
$ws.Range("O217").Value = 'Región de Arica y Parinacota'
$ws.Range("P217").Value = 262
$ws.Range("Q217").Value = 50
$ws.Range("R217").Value = 'Hortaliza'

# New row 218
$ws.Range("A218").Value = 9
$ws.Range("B218").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C218").Value = 'Metropolitana'
$ws.Range("D218").Value = 44845
$ws.Range("E218").Value = 13
$ws.Range("F218").Value = 100112001
$ws.Range("G218").Value = 'Berenjena'
$ws.Range("H218").Value = 'Sin especificar'
$ws.Range("I218").Value = 'Segunda'
$ws.Range("J218").Value = 100
$ws.Range("K218").Value = 8000
$ws.Range("L218").Value = 8000
$ws.Range("M218").Value = 8000
$ws.Range("N218").Value = '$/caja 100 unidades'
$ws.Range("O218").Value = 'Región de Arica y Parinacota'
$ws.Range("P218").Value = 80
$ws.Range("Q218").Value = 100
$ws.Range("R218").Value = 'Hortaliza'
